# nit reported from Pierre
#
# Slide 2 ("Content Placeholder 1", shape id 4), first bullet paragraph:
#   "A single link down result in multiple device topology, ..."
# becomes
#   "A single link down results in multiple device topology, ..."
# with the bold lead-in now split as "A single link " / "down " and a
# new (non-bold) "results " run starting the rest of the sentence.

$p = $ppt.ActivePresentation

# Locate the shape holding the bullet list, searching by its current
# text rather than hard-coding slide/shape indices.
$shp = $null
for ($si = 1; $si -le $p.Slides.Count -and $shp -eq $null; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $cand = $slide.Shapes.Item($shi)
        if ($cand.HasTextFrame -and $cand.TextFrame.HasText) {
            if ($cand.TextFrame.TextRange.Text -like "A single link down result in multiple*") {
                $shp = $cand
                break
            }
        }
    }
}

if ($shp -eq $null) {
    # Fall back to the known location if the text was already edited
    # or could not be matched.
    $shp = $p.Slides.Item(2).Shapes.Item(1)
}

$para1 = $shp.TextFrame.TextRange.Paragraphs(1, 1)

# Fix "result " -> "results " (also splits off the trailing, non-bold
# remainder of the sentence into its own run).
$resultWord = $para1.Characters(20, 7)
if ($resultWord.Text -eq "result ") {
    $resultWord.Text = "results "
}

# Split the leading bold run into "A single link " + "down " so the
# bold formatting ends right after "down ".
$leadIn = $para1.Characters(1, 14)
if ($leadIn.Text -eq "A single link ") {
    $leadIn.Font.Bold = $true
}
